$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '''29.799.30'
$ws.Range('E2').Value = '''  -0.53%  '
$ws.Range('D3').Value = '''1.872.29'
$ws.Range('E3').Value = '''  -0.28%  '
$ws.Range('D4').Value = '''1.000'
$ws.Range('E4').Value = '''  +0.03%  '
$ws.Range('D5').Value = '''0.7321'
$ws.Range('E5').Value = '''  -1.32%  '
$ws.Range('D6').Value = '''241.10'
$ws.Range('E6').Value = '''  -0.45%  '
$ws.Range('D7').Value = '''0.9998'
$ws.Range('E7').Value = '''  +0.05%  '
$ws.Range('E8').Value = '''  -0.58%  '
$ws.Range('D9').Value = '''0.07127'
$ws.Range('E9').Value = '''  -0.52%  '
$ws.Range('E10').Value = '''  -1.59%  '
$ws.Range('D11').Value = '''0.08167'
$ws.Range('E11').Value = '''  -3.17%  '
$ws.Range('B12').Value = '''WrappedEther'
$ws.Range('C12').Value = '''https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D12').Value = '''1.909.06'
$ws.Range('E12').Value = '''  +1.79%  '
$ws.Range('B13').Value = '''Polygon'
$ws.Range('C13').Value = '''https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D13').Value = '''0.7426'
$ws.Range('E13').Value = '''  -1.71%  '
$ws.Range('D14').Value = '''5.343'
$ws.Range('E14').Value = '''  -1.70%  '
$ws.Range('D15').Value = '''92.43'
$ws.Range('D16').Value = '''29.829.37'
$ws.Range('E16').Value = '''  -0.37%  '
$ws.Range('D17').Value = '''6.013'
$ws.Range('E17').Value = '''  -1.40%  '
$ws.Range('D18').Value = '''248.25'
$ws.Range('E18').Value = '''  +1.67%  '
$ws.Range('E19').Value = '''  -1.82%  '
$ws.Range('D20').Value = '''0.000007804'
$ws.Range('E20').Value = '''  -0.43%  '
$ws.Range('D21').Value = '''2.159.78'
$ws.Range('E21').Value = '''  +2.33%  '
$ws.Range('E22').Value = '''  +0.33%  '
$ws.Range('D23').Value = '''1.000'
$ws.Range('E23').Value = '''  +0.05%  '
$ws.Range('D24').Value = '''7.768'
$ws.Range('E24').Value = '''  -2.88%  '
$ws.Range('D25').Value = '''0.1543'
$ws.Range('E25').Value = '''  -1.52%  '
$ws.Range('D26').Value = '''9.201'
$ws.Range('E26').Value = '''  -1.40%  '
$ws.Range('D27').Value = '''163.86'
$ws.Range('E27').Value = '''  -0.43%  '
$ws.Range('D28').Value = '''18.53'
$ws.Range('E28').Value = '''  -0.74%  '
$ws.Range('D29').Value = '''2.019'
$ws.Range('E29').Value = '''  -1.05%  '
$ws.Range('D30').Value = '''1.447'
$ws.Range('E30').Value = '''  -1.88%  '
$ws.Range('D31').Value = '''4.524'
$ws.Range('E31').Value = '''  -1.98%  '
$ws.Range('E32').Value = '''  -0.63%  '
$ws.Range('D33').Value = '''4.187'
$ws.Range('E33').Value = '''  -2.35%  '
$ws.Range('E34').Value = '''  -0.61%  '
$ws.Range('E35').Value = '''  -0.69%  '
$ws.Range('E36').Value = '''  -2.14%  '
$ws.Range('D37').Value = '''1.002'
$ws.Range('E37').Value = '''  -0.13%  '
$ws.Range('D38').Value = '''2.697'
$ws.Range('E38').Value = '''  +0.14%  '
$ws.Range('D39').Value = '''0.01934'
$ws.Range('E39').Value = '''  -1.14%  '
$ws.Range('D40').Value = '''2.732'
$ws.Range('E40').Value = '''  -0.69%  '
$ws.Range('D41').Value = '''0.4462'
$ws.Range('E41').Value = '''  -0.75%  '
$ws.Range('D42').Value = '''5.977'
$ws.Range('E42').Value = '''  -2.48%  '
$ws.Range('D43').Value = '''0.8679'
$ws.Range('E43').Value = '''  +0.38%  '
$ws.Range('D44').Value = '''71.36'
$ws.Range('E44').Value = '''  -1.82%  '
$ws.Range('D45').Value = '''1.044.47'
$ws.Range('E45').Value = '''  -5.87%  '
$ws.Range('B46').Value = '''Quant'
$ws.Range('C46').Value = '''https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D46').Value = '''103.96'
$ws.Range('E46').Value = '''  +0.71%  '
$ws.Range('B47').Value = '''PaxDollar'
$ws.Range('C47').Value = '''https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D47').Value = '''1.0000'
$ws.Range('E47').Value = '''  -0.08%  '
$ws.Range('D48').Value = '''1.818'
$ws.Range('E48').Value = '''  -1.85%  '
$ws.Range('D49').Value = '''7.434'
$ws.Range('E49').Value = '''  -3.66%  '
$ws.Range('D50').Value = '''9.508'
$ws.Range('E50').Value = '''  -0.55%  '
$ws.Range('D51').Value = '''2.047.95'
$ws.Range('E51').Value = '''  +1.76%  '
